# ---------------------------------------------------------------------------
# Adds a new "CONSTANT-GROWTH" LGD model and an "INDEXED" LGD model:
#  - renames the ASSUMPTIONS.lgd_collateral_index field to lgd_index
#  - inserts a new lgd_growth_rate field before it
#  - updates the DICTIONARY sheet to match (new row + updated descriptions)
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. ASSUMPTIONS sheet: insert a new column M (lgd_growth_rate),
#    shifting the old lgd_collateral_index column to N (renamed lgd_index)
# ------------------------------------------------------------------
$assumptions = $wb.Worksheets.Item("ASSUMPTIONS")

$assumptions.Columns("M:M").Insert()
$assumptions.Columns("M:M").ColumnWidth = 17.166666666666668

$assumptions.Range("M1").Value2 = "lgd_growth_rate"
$assumptions.Range("N1").Value2 = "lgd_index"

$assumptions.Range("N2").Select()

# ------------------------------------------------------------------
# 2. DICTIONARY sheet: insert a row describing lgd_growth_rate right
#    after lgd_loss_given_default, rename lgd_collateral_index row to
#    lgd_index, and refresh the lgd_type description.
# ------------------------------------------------------------------
$dictionary = $wb.Worksheets.Item("DICTIONARY")

$dictionary.Rows("14:14").Insert()
$dictionary.Range("A13:C13").Copy()
$dictionary.Range("A14:C14").PasteSpecial(-4122)

$dictionary.Range("A14").Value2 = "lgd_growth_rate"
$dictionary.Range("B14").Value2 = "The constant growth rate to use for the CONSTANT-GROWTH LGD model. The values is expressed as an effective annual rate. "
$dictionary.Range("C14").Value2 = "<float>"
$dictionary.Rows("14:14").RowHeight = 30

$dictionary.Range("A15").Value2 = "lgd_index"

$dictionary.Range("B12").Value2 = "The LGD model to apply. The following values are supported:`n- CONSTANT`n- CONSTANT-GROWTH`n- INDEXED`n- SECURED`n- UNSECURED"
$dictionary.Rows("12:12").RowHeight = 90

$dictionary.Range("B13").Select()
